$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$questions = @(
  "What's the maximum curves I can load in a data file?",
  "I want to add this name as my curve shade name length `"Hydrocarbon bearing zone highlighted`"",
  "What's the maximum number of data points allowed per curve?",
  "Why can't I add 251 curve shades to my log?",
  "I have 20000 modifiers added ty log, why I can't I add anymore?",
  "How many log headers can I add to my log?",
  "How many tadpole definitions am I allowed?",
  "Why can I not add another layout to my log?"
)

$models = @(
  "llama3.2:latest",
  "llama3.2:latest",
  "llama3.2:latest",
  "llama3.2:latest",
  "llama3.2:latest",
  "llama3.2:latest",
  "llama3.2:latest",
  "llama3.2:latest"
)

$responses = @(
  "Based on Document 29, which states: `"Deviation from the above may cause errors during loading data from files.`" and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is:`nYou can load up to 450 curves per data file.",
  "To set the curve shading name, click on the `"Curve`" menu and select `"Shading`". Then, in the `"Shading`" dialog box, enter `"Hydrocarbon bearing zone highlighted`" in the `"Name`" field. Note that the character limit for the curve shade name is 20 characters; ensure your chosen name does not exceed this limit to avoid any errors or issues with the application.",
  "Based on Document 29, which states: `"Deviation from the above may cause errors during loading data from files.`" and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is:`nYou can load an unlimited number of data points per curve.",
  "According to Document 25, the maximum number of curve shades per plot is 250. This means you cannot add more than 250 curve shades to your log.",
  "You cannot add more than 20000 modifiers per plot because of the limit specified on theHometab.",
  "Based on the provided feedback, it seems that the correct answer is not explicitly stated in the documentation. However, it does mention that you can select up to 50 headers for display in a log. Therefore, the optimized answer would be:`n`"The maximum number of headers you can add to your log is 50.`"",
  "You are allowed to create a maximum of 5 new tadpole definition templates.",
  "According to the document, you cannot add another layout to your log because the layout specifies up to 19 layouts per ODF file."
)

$startRow = 255
for ($i = 0; $i -lt $questions.Length; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $questions[$i]
  $ws.Cells.Item($r, 2).Value = $models[$i]
  $ws.Cells.Item($r, 3).Value = $responses[$i]
}

$wb.Save()
